$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the updated cells keep their original text/string storage
# (matching columns D, E, G for data rows 2-51) instead of being
# auto-converted to numbers/percentages by Excel's smart-entry parsing.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "300.77"
$ws.Range("E2").Value = "2.31%"
$ws.Range("G2").Value = "22"
$ws.Range("D3").Value = "32.14"
$ws.Range("E3").Value = "2.83%"
$ws.Range("G3").Value = "22"
$ws.Range("D4").Value = "4.994"
$ws.Range("E4").Value = "0.91%"
$ws.Range("G4").Value = "22"
$ws.Range("D5").Value = "0.07755"
$ws.Range("E5").Value = "5.79%"
$ws.Range("G5").Value = "22"
$ws.Range("D6").Value = "2.287"
$ws.Range("E6").Value = "0.43%"
$ws.Range("G6").Value = "22"
$ws.Range("D7").Value = "7.958"
$ws.Range("E7").Value = "3.28%"
$ws.Range("G7").Value = "22"
$ws.Range("D8").Value = "0.9314"
$ws.Range("E8").Value = "2.58%"
$ws.Range("G8").Value = "22"
$ws.Range("D9").Value = "0.1012"
$ws.Range("E9").Value = "25.00%"
$ws.Range("G9").Value = "22"
$ws.Range("D10").Value = "0.1770"
$ws.Range("E10").Value = "5.10%"
$ws.Range("G10").Value = "22"
$ws.Range("D11").Value = "0.08465"
$ws.Range("E11").Value = "3.52%"
$ws.Range("G11").Value = "22"
$ws.Range("D12").Value = "0.03280"
$ws.Range("E12").Value = "5.63%"
$ws.Range("G12").Value = "22"
$ws.Range("D13").Value = "0.09888"
$ws.Range("E13").Value = "-1.91%"
$ws.Range("G13").Value = "22"
$ws.Range("E14").Value = "-2.07%"
$ws.Range("G14").Value = "22"
$ws.Range("D15").Value = "0.005698"
$ws.Range("E15").Value = "-0.98%"
$ws.Range("G15").Value = "22"
$ws.Range("D16").Value = "3.506"
$ws.Range("E16").Value = "0.67%"
$ws.Range("G16").Value = "22"
$ws.Range("D17").Value = "3.828"
$ws.Range("E17").Value = "2.42%"
$ws.Range("G17").Value = "22"
$ws.Range("D18").Value = "2.197"
$ws.Range("E18").Value = "5.65%"
$ws.Range("G18").Value = "22"
$ws.Range("D19").Value = "0.3367"
$ws.Range("E19").Value = "1.18%"
$ws.Range("G19").Value = "22"
$ws.Range("E20").Value = "3.99%"
$ws.Range("G20").Value = "22"
$ws.Range("D21").Value = "4.369"
$ws.Range("E21").Value = "9.71%"
$ws.Range("G21").Value = "22"
$ws.Range("D22").Value = "0.2088"
$ws.Range("E22").Value = "-0.73%"
$ws.Range("G22").Value = "22"
$ws.Range("D23").Value = "0.04584"
$ws.Range("E23").Value = "1.21%"
$ws.Range("G23").Value = "22"
$ws.Range("E24").Value = "0.47%"
$ws.Range("G24").Value = "22"
$ws.Range("D25").Value = "0.004370"
$ws.Range("E25").Value = "0.64%"
$ws.Range("G25").Value = "22"
$ws.Range("E26").Value = "-0.63%"
$ws.Range("G26").Value = "22"
$ws.Range("D27").Value = "0.0003374"
$ws.Range("G27").Value = "22"
$ws.Range("G28").Value = "22"
$ws.Range("G29").Value = "22"
$ws.Range("G30").Value = "22"
$ws.Range("G31").Value = "22"
$ws.Range("G32").Value = "22"
$ws.Range("G33").Value = "22"
$ws.Range("G34").Value = "22"
$ws.Range("G35").Value = "22"
$ws.Range("G36").Value = "22"
$ws.Range("G37").Value = "22"
$ws.Range("G38").Value = "22"
$ws.Range("D39").Value = "0.01699"
$ws.Range("E39").Value = "5.89%"
$ws.Range("G39").Value = "22"
$ws.Range("D40").Value = "0.04718"
$ws.Range("E40").Value = "6.35%"
$ws.Range("G40").Value = "22"
$ws.Range("D41").Value = "0.007719"
$ws.Range("E41").Value = "4.82%"
$ws.Range("G41").Value = "22"
$ws.Range("D42").Value = "0.009773"
$ws.Range("E42").Value = "8.26%"
$ws.Range("G42").Value = "22"
$ws.Range("D43").Value = "0.1392"
$ws.Range("E43").Value = "4.76%"
$ws.Range("G43").Value = "22"
$ws.Range("D44").Value = "0.002089"
$ws.Range("E44").Value = "7.02%"
$ws.Range("G44").Value = "22"
$ws.Range("D45").Value = "0.009663"
$ws.Range("E45").Value = "1.51%"
$ws.Range("G45").Value = "22"
$ws.Range("D46").Value = "0.00006080"
$ws.Range("E46").Value = "2.16%"
$ws.Range("G46").Value = "22"
$ws.Range("E47").Value = "-0.74%"
$ws.Range("G47").Value = "22"
$ws.Range("D48").Value = "2.655"
$ws.Range("E48").Value = "18.45%"
$ws.Range("G48").Value = "22"
$ws.Range("D49").Value = "0.001988"
$ws.Range("E49").Value = "-31.42%"
$ws.Range("G49").Value = "22"
$ws.Range("D50").Value = "0.00002088"
$ws.Range("E50").Value = "-0.74%"
$ws.Range("G50").Value = "22"
$ws.Range("D51").Value = "0.0001988"
$ws.Range("E51").Value = "-0.74%"
$ws.Range("G51").Value = "22"
